$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 7000
$ws.Range("I16").Value = 7000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -6770
$ws.Range("N16").ClearContents()

# Row 21
$ws.Range("H21").Value = 17000.143
$ws.Range("I21").Value = 7067
$ws.Range("J21").Value = 24450
$ws.Range("K21").Value = 7067
$ws.Range("L21").Value = 24450
$ws.Range("M21").Value = -6599
$ws.Range("N21").Value = -25386

# Row 23
$ws.Range("H23").Value = 17000.143
$ws.Range("I23").Value = 7067
$ws.Range("J23").Value = 24450
$ws.Range("K23").Value = 7067
$ws.Range("L23").Value = 24450
$ws.Range("M23").Value = -6833
$ws.Range("N23").Value = -24918

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()

# Row 32
$ws.Range("H32").Value = 748.6667
$ws.Range("I32").Value = 366.66666
$ws.Range("J32").Value = 2276.6667
$ws.Range("K32").Value = 366.66666
$ws.Range("L32").Value = 2276.6667
$ws.Range("M32").Value = -40.66665999999998
$ws.Range("N32").Value = -2928.6667

# Row 33
$ws.Range("H33").Value = 144.8125
$ws.Range("I33").Value = 61.727272
$ws.Range("J33").Value = 327.6
$ws.Range("K33").Value = 61.727272
$ws.Range("L33").Value = 327.6
$ws.Range("M33").Value = 167.272728
$ws.Range("N33").Value = -785.6

# Row 39
$ws.Range("H39").Value = 168.13333
$ws.Range("I39").Value = 117.84615
$ws.Range("J39").Value = 495
$ws.Range("K39").Value = 353.53845
$ws.Range("L39").Value = 1485
$ws.Range("M39").Value = -57.53845000000001
$ws.Range("N39").Value = -2077

# Row 40
$ws.Range("H40").Value = 3160
$ws.Range("I40").Value = 3160
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3160
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2985
$ws.Range("N40").ClearContents()

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 51
$ws.Range("H51").Value = 2546.9
$ws.Range("I51").Value = 1999.75
$ws.Range("J51").Value = 2911.6667
$ws.Range("K51").Value = 1999.75
$ws.Range("L51").Value = 2911.6667
$ws.Range("M51").Value = -1515.75
$ws.Range("N51").Value = -3879.6667

# Row 52
$ws.Range("H52").Value = 3000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 3000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 9000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -9320

# Row 58
$ws.Range("H58").Value = 2859.4285
$ws.Range("I58").Value = 105
$ws.Range("J58").Value = 3610.6365
$ws.Range("K58").Value = 315
$ws.Range("L58").Value = 10831.9095
$ws.Range("M58").Value = -165
$ws.Range("N58").Value = -11131.9095

# Row 132
$ws.Range("H132").Value = 2029.7742
$ws.Range("I132").Value = 2073.2932
$ws.Range("J132").Value = 1398.75
$ws.Range("K132").Value = 6219.8796
$ws.Range("L132").Value = 4196.25
$ws.Range("M132").Value = -3689.8796
$ws.Range("N132").Value = -9256.25

# Row 138
$ws.Range("H138").Value = 3779.5325
$ws.Range("I138").Value = 1746.0667
$ws.Range("J138").Value = 4271.5
$ws.Range("K138").Value = 5238.2001
$ws.Range("L138").Value = 12814.5
$ws.Range("M138").Value = -98.20010000000002
$ws.Range("N138").Value = -23094.5


$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 308839.9
$ws.Range("I6").Value = 3000000
$ws.Range("J6").Value = 9822.111000000001
$ws.Range("K6").Value = 3000000
$ws.Range("L6").Value = 9822.111000000001
$ws.Range("M6").Value = -2999827
$ws.Range("N6").Value = -10168.111

# Row 26
$ws.Range("H26").Value = 4341.1665
$ws.Range("I26").Value = 1209.4
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 1209.4
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -879.4000000000001
$ws.Range("N26").Value = -20660

# Row 54
$ws.Range("H54").Value = 20316.334
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 20316.334
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 20316.334
$ws.Range("N54").Value = -21854.334


$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 7950.3335
$ws.Range("I7").Value = 1975
$ws.Range("J7").Value = 10938
$ws.Range("K7").Value = 1975
$ws.Range("L7").Value = 10938
$ws.Range("M7").Value = -1862
$ws.Range("N7").Value = -11164


$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 14546.667
$ws.Range("I2").Value = 640
$ws.Range("J2").Value = 21500
$ws.Range("K2").Value = 640
$ws.Range("L2").Value = 21500
$ws.Range("M2").Value = -527
$ws.Range("N2").Value = -21726

# Row 6
$ws.Range("H6").Value = 1763
$ws.Range("I6").Value = 953.6667
$ws.Range("J6").Value = 2370
$ws.Range("K6").Value = 953.6667
$ws.Range("L6").Value = 2370
$ws.Range("M6").Value = -840.6667
$ws.Range("N6").Value = -2596

# Row 12
$ws.Range("H12").Value = 3475.75
$ws.Range("I12").Value = 766.6667
$ws.Range("J12").Value = 5101.2
$ws.Range("K12").Value = 766.6667
$ws.Range("L12").Value = 5101.2
$ws.Range("M12").Value = -596.6667
$ws.Range("N12").Value = -5441.2

# Row 19
$ws.Range("H19").Value = 175.21739
$ws.Range("I19").Value = 219.54546
$ws.Range("J19").Value = 134.58333
$ws.Range("K19").Value = 219.54546
$ws.Range("L19").Value = 134.58333
$ws.Range("M19").Value = -49.54545999999999
$ws.Range("N19").Value = -474.58333

# Row 24
$ws.Range("H24").Value = 175.21739
$ws.Range("I24").Value = 219.54546
$ws.Range("J24").Value = 134.58333
$ws.Range("K24").Value = 219.54546
$ws.Range("L24").Value = 134.58333
$ws.Range("M24").Value = -49.54545999999999
$ws.Range("N24").Value = -474.58333

# Row 33
$ws.Range("H33").Value = 33005.832
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 33005.832
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 33005.832
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -33763.832

# Row 36
$ws.Range("H36").Value = 13265.3
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 13265.3
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 13265.3
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -14041.3

# Row 40
$ws.Range("H40").Value = 13265.3
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 13265.3
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 13265.3
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -13585.3

# Row 50
$ws.Range("H50").Value = 20000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250

# Row 51
$ws.Range("H51").Value = 12000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 12000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 12000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -13472

# Row 59
$ws.Range("H59").Value = 12845.211
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 12845.211
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 12845.211
$ws.Range("N59").Value = -15135.211

# Row 60
$ws.Range("H60").Value = 9933.333000000001
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 9933.333000000001
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 9933.333000000001
$ws.Range("N60").Value = -10955.333

# Row 61
$ws.Range("H61").Value = 12000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 12000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 12000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -12696

# Row 68
$ws.Range("H68").Value = 17422.666
$ws.Range("I68").Value = 8634
$ws.Range("J68").Value = 35000
$ws.Range("K68").Value = 8634
$ws.Range("L68").Value = 35000
$ws.Range("M68").Value = -7885
$ws.Range("N68").Value = -36498

# Row 70
$ws.Range("H70").Value = 21975
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 21975
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 21975
$ws.Range("N70").Value = -22605

# Row 71
$ws.Range("H71").Value = 17422.666
$ws.Range("I71").Value = 8634
$ws.Range("J71").Value = 35000
$ws.Range("K71").Value = 25902
$ws.Range("L71").Value = 105000
$ws.Range("M71").Value = -22158
$ws.Range("N71").Value = -112488

# Row 73
$ws.Range("H73").Value = 21975
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 21975
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 21975
$ws.Range("N73").Value = -24159


$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 6463.375
$ws.Range("I9").Value = 421.4
$ws.Range("J9").Value = 16533.334
$ws.Range("K9").Value = 421.4
$ws.Range("L9").Value = 16533.334
$ws.Range("M9").Value = -251.4
$ws.Range("N9").Value = -16873.334

# Row 31
$ws.Range("H31").Value = 3578.6
$ws.Range("I31").Value = 3223.25
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 3223.25
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2931.25
$ws.Range("N31").Value = -5584

# Row 37
$ws.Range("H37").Value = 3578.6
$ws.Range("I37").Value = 3223.25
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 3223.25
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -2946.25
$ws.Range("N37").Value = -5554


$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 7576.7144
$ws.Range("I9").Value = 305
$ws.Range("J9").Value = 13030.5
$ws.Range("K9").Value = 305
$ws.Range("L9").Value = 13030.5
$ws.Range("M9").Value = -81
$ws.Range("N9").Value = -13478.5

# Row 30
$ws.Range("H30").Value = 13373.714
$ws.Range("I30").Value = 904
$ws.Range("J30").Value = 30000
$ws.Range("K30").Value = 904
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = -796
$ws.Range("N30").Value = -30216

# Row 35
$ws.Range("H35").Value = 15719.857
$ws.Range("I35").Value = 1259.75
$ws.Range("J35").Value = 35000
$ws.Range("K35").Value = 1259.75
$ws.Range("L35").Value = 35000
$ws.Range("M35").Value = -923.75
$ws.Range("N35").Value = -35672

